$d = $word.ActiveDocument
$d.Content.Find.Execute("584×5=", $true, $false, $false, $false, $false, $true, 1, $false, "835×9=", 2)
$d.Content.Find.Execute("299×9=", $true, $false, $false, $false, $false, $true, 1, $false, "822×5=", 2)
$d.Content.Find.Execute("488×7=", $true, $false, $false, $false, $false, $true, 1, $false, "340×6=", 2)
$d.Content.Find.Execute("293×5=", $true, $false, $false, $false, $false, $true, 1, $false, "216×6=", 2)
$d.Content.Find.Execute("208×7=", $true, $false, $false, $false, $false, $true, 1, $false, "369×6=", 2)
$d.Content.Find.Execute("365×9=", $true, $false, $false, $false, $false, $true, 1, $false, "146×8=", 2)
$d.Content.Find.Execute("236×9=", $true, $false, $false, $false, $false, $true, 1, $false, "749×9=", 2)
$d.Content.Find.Execute("529×4=", $true, $false, $false, $false, $false, $true, 1, $false, "971×2=", 2)
$d.Content.Find.Execute("398×7=", $true, $false, $false, $false, $false, $true, 1, $false, "556×6=", 2)
$d.Content.Find.Execute("605×7=", $true, $false, $false, $false, $false, $true, 1, $false, "474×6=", 2)
$d.Content.Find.Execute("740×5=", $true, $false, $false, $false, $false, $true, 1, $false, "751×3=", 2)
$d.Content.Find.Execute("182×4=", $true, $false, $false, $false, $false, $true, 1, $false, "411×3=", 2)
$d.Content.Find.Execute("163×6=", $true, $false, $false, $false, $false, $true, 1, $false, "386×3=", 2)
$d.Content.Find.Execute("898×6=", $true, $false, $false, $false, $false, $true, 1, $false, "248×9=", 2)
$d.Content.Find.Execute("361×5=", $true, $false, $false, $false, $false, $true, 1, $false, "357×6=", 2)
$d.Content.Find.Execute("161×9=", $true, $false, $false, $false, $false, $true, 1, $false, "601×6=", 2)
$d.Content.Find.Execute("666×4=", $true, $false, $false, $false, $false, $true, 1, $false, "823×3=", 2)
$d.Content.Find.Execute("864×4=", $true, $false, $false, $false, $false, $true, 1, $false, "975×9=", 2)
$d.Content.Find.Execute("795×6=", $true, $false, $false, $false, $false, $true, 1, $false, "764×3=", 2)
$d.Content.Find.Execute("558×8=", $true, $false, $false, $false, $false, $true, 1, $false, "781×3=", 2)
$d.Content.Find.Execute("873×9=", $true, $false, $false, $false, $false, $true, 1, $false, "439×2=", 2)
$d.Content.Find.Execute("270×2=", $true, $false, $false, $false, $false, $true, 1, $false, "678×4=", 2)
$d.Content.Find.Execute("445×6=", $true, $false, $false, $false, $false, $true, 1, $false, "120×8=", 2)
$d.Content.Find.Execute("432×8=", $true, $false, $false, $false, $false, $true, 1, $false, "602×3=", 2)
$d.Content.Find.Execute("916×8=", $true, $false, $false, $false, $false, $true, 1, $false, "713×6=", 2)
